# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Arándano (blue)" at Vega Central
# Mapocho de Santiago, pushing the existing row 154 (and everything below it)
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 154; rows 154..179 shift to 155..180.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(154, 1).Value  = 9
$ws.Cells.Item(154, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(154, 3).Value  = "Metropolitana"
$ws.Cells.Item(154, 4).Value  = 44617
$ws.Cells.Item(154, 5).Value  = 13
$ws.Cells.Item(154, 6).Value  = "Fruta"
$ws.Cells.Item(154, 7).Value  = 100101
$ws.Cells.Item(154, 8).Value  = "Berries"
$ws.Cells.Item(154, 9).Value  = 100101001
$ws.Cells.Item(154, 10).Value = "Arándano (blue)"
$ws.Cells.Item(154, 11).Value = "Sin especificar"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 450
$ws.Cells.Item(154, 14).Value = 3800
$ws.Cells.Item(154, 15).Value = 3800
$ws.Cells.Item(154, 16).Value = 3800
$ws.Cells.Item(154, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(154, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(154, 19).Value = 1900
$ws.Cells.Item(154, 20).Value = 2
